$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# New header cells: Wins / Losses / Ties (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
Set-HeaderStyle($ws.Range("AD1"))

$ws.Range("AE1").Value = "Losses"
Set-HeaderStyle($ws.Range("AE1"))

$ws.Range("AF1").Value = "Ties"
Set-HeaderStyle($ws.Range("AF1"))

# Fill in the season record (Wins=64, Losses=98, Ties=0) for every player row
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 0
}
